# Remove the leading pandoc-style title block that the download script now
# produces separately (heading "On Pilgrimage - January 1949" + the bold
# "By Dorothy Day" byline paragraph, along with the bookmark that wrapped
# the heading). The rest of the article body is left untouched.

$d = $word.ActiveDocument

# Paragraph 1 is the Heading1 "On Pilgrimage - January 1949" paragraph and
# paragraph 2 is the bold "By Dorothy Day" byline paragraph. Deleting
# paragraph 1 twice removes both, since each delete shifts the following
# paragraph up into slot 1.
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(1).Range.Delete()
